$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "task last day" column (D) for the existing rows 2-23 ---
# All rows move from 2023-12-11 (45271) to 2023-12-22 (45282)
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 4).Value = 45282
}

# --- Update the remaining/ratio figures that changed on specific rows ---
# Row 3 (M-MN 管理)
$ws.Range("F3").Value = 67
$ws.Range("H3").Value = 0.8375

# Row 4 (M-BT 出張)
$ws.Range("F4").Value = 28
$ws.Range("H4").Value = 0.35

# Row 5 (M-MT ミーティング)
$ws.Range("F5").Value = 35.25
$ws.Range("H5").Value = 0.2203125

# Row 12 (PP-EL-MP メイン基板)
$ws.Range("F12").Value = -49
$ws.Range("H12").Value = -3.0625

# Row 14 (PP-EL-PP 電源)
$ws.Range("F14").Value = -3
$ws.Range("H14").Value = -0.375

# --- Append new rows 24-26 ---

# Row 24: ES2-PL-PL-EL エレキ
$ws.Cells.Item(24, 1).Value = "ES2-PL-PL-EL エレキ"
$ws.Cells.Item(24, 2).Value = 45272
$ws.Cells.Item(24, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 3).Value = 45273
$ws.Cells.Item(24, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 4).Value = 45282
$ws.Cells.Item(24, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 5).Value = "Expired"
$ws.Cells.Item(24, 6).Value = -1.5
$ws.Cells.Item(24, 7).Value = 16
$ws.Cells.Item(24, 8).Value = -0.09375

# Row 25: ES2-PL-PL-RV レビュー
$ws.Cells.Item(25, 1).Value = "ES2-PL-PL-RV レビュー"
$ws.Cells.Item(25, 2).Value = 45272
$ws.Cells.Item(25, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(25, 3).Value = 45273
$ws.Cells.Item(25, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(25, 4).Value = 45282
$ws.Cells.Item(25, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(25, 5).Value = "Expired"
$ws.Cells.Item(25, 6).Value = 9
$ws.Cells.Item(25, 7).Value = 12
$ws.Cells.Item(25, 8).Value = 0.75

# Row 26: ES2-PL-CD-CD 構想設計資料
$ws.Cells.Item(26, 1).Value = "ES2-PL-CD-CD 構想設計資料"
$ws.Cells.Item(26, 2).Value = 45280
$ws.Cells.Item(26, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(26, 3).Value = 45281
$ws.Cells.Item(26, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(26, 4).Value = 45282
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(26, 5).Value = "Expired"
$ws.Cells.Item(26, 6).Value = 12.5
$ws.Cells.Item(26, 7).Value = 16
$ws.Cells.Item(26, 8).Value = 0.78125
